$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.204.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.783.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.57%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.99"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.550"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.15"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0656"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.040.27"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.16"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.790.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.210.93"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.20"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.98"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0739"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.33"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.38%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.35"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.58%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.58"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.438.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.623"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.79"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.887"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0507"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.939.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.18"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.48"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.26"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.60%  "
